# Applies crypto price/volume updates per the diff (commit: "Updated symbol list on Fri Feb 17 11:25:01 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume/1h) hold numeric-looking text (e.g. "310.12", "-3.39%").
# Force text format first so Excel stores them as literal strings, matching the source data.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "310.12"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-3.39%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "50.75"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "4.90%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.169"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.74%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07793"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-3.72%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.498"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-2.13%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.349"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "11.64%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.567"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-4.87%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1213"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-6.63%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1984"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.54%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.04798"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "4.15%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09458"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.58%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1043"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.54%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001265"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-5.20%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005779"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.60%"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2,015.36%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.329"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.23%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.28%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3474"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.57%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.067"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.07%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1368"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.78%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.3093"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.03%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04154"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.34%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001271"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-2.63%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.003948"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-7.03%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001349"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.11%"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-4.14%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06012"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "4.84%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01099"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "74.39%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007910"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.08%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1423"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.22%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008371"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "8.53%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008332"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "2.85%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3379"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "5.79%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007266"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "5.33%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.11%"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002618"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-34.59%"
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05317"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-19.71%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.11%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.11%"
